$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing data on the sheet
$ws.Cells.Clear()

# Set header values
$ws.Range("A1").Value = "Username "
$ws.Range("B1").Value = "Password"

# Make the header row bold
$ws.Range("A1:B1").Font.Bold = $true

# Update the selection to match the new layout
$ws.Range("A1:B1").Select()
